# Insert a new first column on Sheet1, shifting the existing A:J data to
# B:K, then label each shifted data row with the source column letter it
# originally came from (row 1 gets the header "col").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns.Item(1).Insert()

$ws.Range("A1").Value = "col"
$ws.Range("A2").Value = "C"
$ws.Range("A3").Value = "D"
$ws.Range("A4").Value = "E"
$ws.Range("A5").Value = "F"
$ws.Range("A6").Value = "G"
$ws.Range("A7").Value = "I"
$ws.Range("A8").Value = "J"

$ws.Range("B5").Select()
